$d = $word.ActiveDocument

$replacements = @(
    @("362×9=3258", "838×9=7542"),
    @("250×8=2000", "769×5=3845"),
    @("639×7=4473", "141×3=423"),
    @("364×6=2184", "879×7=6153"),
    @("815×8=6520", "293×6=1758"),
    @("823×8=6584", "635×2=1270"),
    @("326×8=2608", "173×2=346"),
    @("297×4=1188", "768×5=3840"),
    @("722×8=5776", "348×7=2436"),
    @("788×9=7092", "834×4=3336"),
    @("754×3=2262", "664×7=4648"),
    @("120×6=720",  "167×6=1002"),
    @("683×7=4781", "434×6=2604"),
    @("860×3=2580", "350×3=1050"),
    @("985×2=1970", "615×4=2460"),
    @("169×7=1183", "252×4=1008"),
    @("395×8=3160", "566×9=5094"),
    @("765×5=3825", "683×2=1366"),
    @("607×3=1821", "653×3=1959"),
    @("784×7=5488", "946×4=3784"),
    @("180×9=1620", "109×9=981"),
    @("963×7=6741", "746×2=1492"),
    @("765×2=1530", "860×8=6880"),
    @("464×3=1392", "943×2=1886"),
    @("479×4=1916", "555×3=1665")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
